$wb = $excel.ActiveWorkbook

# "Hoja1" already exists; add a new sheet right after it and name it "Hoja2"
$hoja1 = $wb.Worksheets.Item(1)
$hoja2 = $wb.Worksheets.Add([System.Type]::Missing, $hoja1)
$hoja2.Name = "Hoja2"

# Row 3: two 3x1 blocks (3,4,5) side by side, plus a lone 9
$hoja2.Range("C3").Value = 3
$hoja2.Range("D3").Value = 4
$hoja2.Range("E3").Value = 5
$hoja2.Range("G3").Value = 3
$hoja2.Range("H3").Value = 4
$hoja2.Range("I3").Value = 5
$hoja2.Range("K3").Value = 9

# Row 4: same two 3x1 blocks
$hoja2.Range("C4").Value = 3
$hoja2.Range("D4").Value = 4
$hoja2.Range("E4").Value = 5
$hoja2.Range("G4").Value = 3
$hoja2.Range("H4").Value = 4
$hoja2.Range("I4").Value = 5

# Row 5: same two 3x1 blocks
$hoja2.Range("C5").Value = 3
$hoja2.Range("D5").Value = 4
$hoja2.Range("E5").Value = 5
$hoja2.Range("G5").Value = 3
$hoja2.Range("H5").Value = 4
$hoja2.Range("I5").Value = 5

# Leave the selection/active cell on Hoja2 at I7, matching the saved view
$hoja2.Range("I7").Select() | Out-Null
